$wb = $excel.ActiveWorkbook

# 1. Rename the "Include" worksheet tab
$wsInclude = $wb.Worksheets.Item("Include from ActRelationshipT")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata worksheet
$ws = $wb.Worksheets.Item("Metadata")

# Bump the Version and Date metadata values
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row just before "Description"
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
